# Generate Report for Handoff
#
# This reflects a re-run of the localization status report generation for
# the zh-cn and de-de handoff sheets: the "370ff961", "8d098eb0",
# "a5021f51" and "ba5a2ccf" source files (rows 4-7) just had their
# handoff packages (re)generated, so their "Latest Handoff Datetime"
# stamps move forward and - for zh-cn - the Priority goes from the
# (stale) "low" placeholder to the real "ht" value already used by the
# sibling rows.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 4-7 (Priority + Latest Handoff Datetime) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-08-26 00:31:16"
}

# --- de-de sheet: rows 4-7 (Latest Handoff Datetime) ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $wsDe.Cells.Item($r, 8).Value = "2016-08-26 00:31:21"
}

# --- Overview sheet: rows 4-7 (Latest HO Xliff Generate Date mirrors the
#     de-de handoff stamp, which is the newer of the two locales) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-26 00:31:21"
}
